$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns B-E on this sheet hold text data (coin name, link, price,
# and percentage change), even when a value looks numeric (e.g. "412.35" or
# "0.740"). Force text formatting before assigning so Excel does not silently
# convert the string into a floating point number (which would lose exact
# formatting such as trailing zeros, or introduce floating-point noise).
$cells = @(
    "D2"
    "E2"
    "D3"
    "E3"
    "E4"
    "D5"
    "E5"
    "D6"
    "E6"
    "E7"
    "E8"
    "D9"
    "E9"
    "E10"
    "E11"
    "D12"
    "E12"
    "D13"
    "E13"
    "D14"
    "E14"
    "E15"
    "D16"
    "E16"
    "D17"
    "E17"
    "D18"
    "E18"
    "E19"
    "D20"
    "E20"
    "D21"
    "E21"
    "D22"
    "E22"
    "E23"
    "D24"
    "E24"
    "E25"
    "D26"
    "E26"
    "D27"
    "E27"
    "D28"
    "E28"
    "D29"
    "E29"
    "E30"
    "E31"
    "E32"
    "E33"
    "D34"
    "E34"
    "D35"
    "E35"
    "E36"
    "D37"
    "E37"
    "D38"
    "E38"
    "D39"
    "E39"
    "E40"
    "D41"
    "E41"
    "E42"
    "B43"
    "C43"
    "D43"
    "E43"
    "B44"
    "C44"
    "D44"
    "E44"
    "D45"
    "E45"
    "D46"
    "E46"
    "D47"
    "E47"
    "D48"
    "E48"
    "D49"
    "E49"
    "D50"
    "E50"
    "D51"
    "E51"
)
foreach ($cellref in $cells) {
    $ws.Range($cellref).NumberFormat = "@"
}

$ws.Range("D2").Value = '62.056.93'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '3.440.92'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '412.35'
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("D6").Value = '130.02'
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("E7").Value = '  +1.39%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.740'
$ws.Range("E9").Value = '  -2.47%  '
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").Value = '0.0000223'
$ws.Range("E12").Value = '  +11.34%  '
$ws.Range("D13").Value = '9.37'
$ws.Range("E13").Value = '  +4.34%  '
$ws.Range("D14").Value = '3.988.69'
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("D16").Value = '21.29'
$ws.Range("E16").Value = '  +2.69%  '
$ws.Range("D17").Value = '3.444.35'
$ws.Range("E17").Value = '  +0.32%  '
$ws.Range("D18").Value = '12.71'
$ws.Range("E18").Value = '  +1.55%  '
$ws.Range("E19").Value = '  +1.48%  '
$ws.Range("D20").Value = '62.136.54'
$ws.Range("E20").Value = '  -0.24%  '
$ws.Range("D21").Value = '494.17'
$ws.Range("E21").Value = '  +22.40%  '
$ws.Range("D22").Value = '93.05'
$ws.Range("E22").Value = '  +3.04%  '
$ws.Range("E23").Value = '  +3.38%  '
$ws.Range("D24").Value = '13.64'
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("E25").Value = '  +4.78%  '
$ws.Range("D26").Value = '35.15'
$ws.Range("E26").Value = '  +4.42%  '
$ws.Range("D27").Value = '9.21'
$ws.Range("E27").Value = '  +6.18%  '
$ws.Range("D28").Value = '4.81'
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").Value = '7.62'
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("E30").Value = '  +2.24%  '
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("E32").Value = '  -2.15%  '
$ws.Range("E33").Value = '  -2.06%  '
$ws.Range("D34").Value = '42.23'
$ws.Range("E34").Value = '  -4.17%  '
$ws.Range("D35").Value = '59.70'
$ws.Range("E35").Value = '  +13.74%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '0.0500'
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("D38").Value = '3.50'
$ws.Range("E38").Value = '  +2.42%  '
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("E40").Value = '  +3.97%  '
$ws.Range("D41").Value = '151.31'
$ws.Range("E41").Value = '  +7.63%  '
$ws.Range("E42").Value = '  +8.06%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '2.97'
$ws.Range("E43").Value = '  +1.54%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = '2.73'
$ws.Range("E44").Value = '  +12.80%  '
$ws.Range("D45").Value = '0.319'
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("D46").Value = '4.33'
$ws.Range("E46").Value = '  +6.12%  '
$ws.Range("D47").Value = '2.38'
$ws.Range("E47").Value = '  +22.80%  '
$ws.Range("D48").Value = '16.68'
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("D49").Value = '121.03'
$ws.Range("E49").Value = '  +24.97%  '
$ws.Range("D50").Value = '23.18'
$ws.Range("E50").Value = '  +4.74%  '
$ws.Range("D51").Value = '0.147'
$ws.Range("E51").Value = '  +16.79%  '
